# "Removed the Teresa case as an example"
#
# The "EDA" sheet already has one worked example in F2:F5 that flags names
# whose space-separated words all start with the same (case-sensitive)
# first letter - that formula still matches "Teresa" (a single-word name)
# because a one-element array trivially equals itself.
#
# This edit appends two more worked examples further down the same sheet:
#   - F8  : a broken attempt at refactoring the helper into its own LAMBDA
#           (calls the bare name `x` instead of the LAMBDA's own parameter),
#           which raises #VALUE!, annotated with a note in G8.
#   - F11 : the corrected version, with an extra COLUMNS(z)>1 guard so a
#           single-word name like "Teresa" no longer qualifies, annotated
#           with a note in J11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Broken example: bare `x` instead of `_xlpm.x` inside MAP's LAMBDA -> #VALUE!
$ws.Range("F8").FormulaArray = '=LET(fx, LAMBDA(x, LET(z, LEFT(TEXTSPLIT(x," "),1),IFERROR(AND(z=INDEX(z,,1)),FALSE))), FILTER(A2:A11,MAP(A2:A11,fx(x))))'
$ws.Range("G8").Value = "Can't just use a bare function. It has to be in a lambda."

# Fixed example: adds COLUMNS(z)>1 so single-name rows (e.g. "Teresa") drop out
$ws.Range("F11:F13").FormulaArray = '=LET(fx, LAMBDA(x, LET(z, LEFT(TEXTSPLIT(x," "),1),IFERROR(AND(z=INDEX(z,,1),COLUMNS(z)>1),FALSE))), FILTER(A2:A11,MAP(A2:A11,LAMBDA(x,fx(x)))))'
$ws.Range("J11").Value = "Simple mod to remove single name cases"

$ws.Range("J15").Select()
